# docs/plan.xlsx update
#  - "план" sheet: drop the "В идеале..." comment, shift the remaining
#    row comments from column G into column F and delete the now-empty
#    column G, clear the leftover orange highlight on F3:F4, and tweak
#    the "Разработать подробную схему кода" task text.
#  - "вопросы" sheet: add a new question/answer row and make it the
#    active tab/selection when the workbook is saved.

$wb = $excel.ActiveWorkbook

# ---- "план" sheet (tab 1) ----------------------------------------------
$plan = $wb.Worksheets.Item(1)

# Carry the row comments that live in column G over to column F before
# the column shift; G3's comment is simply dropped (no longer needed).
$commentRow5  = $plan.Range("G5").Value()
$commentRow11 = $plan.Range("G11").Value()
$commentRow12 = $plan.Range("G12").Value()

$plan.Range("F5").Value  = $commentRow5
$plan.Range("F11").Value = $commentRow11
$plan.Range("F12").Value = $commentRow12

$plan.Range("G3").ClearContents()
$plan.Range("G5").ClearContents()
$plan.Range("G11").ClearContents()
$plan.Range("G12").ClearContents()

# Column G is now completely empty - remove it outright.
$plan.Columns("G").Delete()

# Remove the leftover orange "todo" highlight from F3:F4.
$plan.Range("F3:F4").ClearFormats()

# Clarify which schema the task is about.
$plan.Range("B11").Value = "Разработать подробную схему кода Для API, интерфейсов и классов, чтобы всё было унифицировано и дружило друг с другом в любой плоскости"

$plan.Range("B12").Select()

# ---- "вопросы" sheet (tab 5) -------------------------------------------
$questions = $wb.Worksheets.Item(5)

$questions.Range("A4").Value = 3
$questions.Range("B4").Value = "Только сейчас почему-то задалась вопросом как будет выглядеть магазин одежды например(магазин ассетов), мб придётся создать страничку в браузере и в ней можно будет листать и выбирать покупку? Мб как-то в интерфейсе телеграма можно было бы это сделать, у него есть апи для приложений магазинов и оплаты, но что тогда в случае с дискордом делать, потому что там нет такого функционала"

$questions.Range("G11").Select()
$questions.Activate()

$wb.Save()
